$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 'Bitcoin'
$ws.Range("C2").Value = 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'
$ws.Range("D2").Value = '24.425.19'
$ws.Range("E2").Value = '  +8.89%  '

$ws.Range("B3").Value = 'Ethereum'
$ws.Range("C3").Value = 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'
$ws.Range("D3").Value = '1.679.45'
$ws.Range("E3").Value = '  +4.79%  '

$ws.Range("B4").Value = 'TetherUSD'
$ws.Range("C4").Value = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '1.005'
$c.Style = "Normal"
$ws.Range("E4").Value = '  +0.10%  '

$ws.Range("B5").Value = 'USDC'
$ws.Range("C5").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '0.9984'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +0.69%  '

$ws.Range("B6").Value = 'BNB'
$ws.Range("C6").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '306.07'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +0.65%  '

$ws.Range("B7").Value = 'XRP'
$ws.Range("C7").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.3715'
$c.Style = "Normal"
$ws.Range("E7").Value = '  +0.83%  '

$ws.Range("B8").Value = 'Cardano'
$ws.Range("C8").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.3445'
$c.Style = "Normal"
$ws.Range("E8").Value = '  +1.15%  '

$ws.Range("B9").Value = 'OKB'
$ws.Range("C9").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '48.16'
$c.Style = "Normal"
$ws.Range("E9").Value = '  +12.67%  '

$ws.Range("B10").Value = 'Polygon'
$ws.Range("C10").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '1.183'
$c.Style = "Normal"
$ws.Range("E10").Value = '  +3.67%  '

$ws.Range("B11").Value = 'Dogecoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.07278'
$c.Style = "Normal"
$ws.Range("E11").Value = '  +3.22%  '

$ws.Range("B12").Value = 'BinanceUSD'
$ws.Range("C12").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '1.002'
$c.Style = "Normal"
$ws.Range("E12").Value = '  +0.17%  '

$ws.Range("B13").Value = 'Solana'
$ws.Range("C13").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '20.43'
$c.Style = "Normal"
$ws.Range("E13").Value = '  +3.81%  '

$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '6.137'
$c.Style = "Normal"
$ws.Range("E14").Value = '  +3.32%  '

$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '6.754'
$c.Style = "Normal"
$ws.Range("E15").Value = '  +1.87%  '

$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").Value = '1.678.65'
$ws.Range("E16").Value = '  +4.76%  '

$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '0.00001109'
$c.Style = "Normal"
$ws.Range("E17").Value = '  +2.28%  '

$ws.Range("B18").Value = 'Dai'
$ws.Range("C18").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '0.9985'
$c.Style = "Normal"
$ws.Range("E18").Value = '  +0.67%  '

$ws.Range("B19").Value = 'TRON'
$ws.Range("C19").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '0.06720'
$c.Style = "Normal"
$ws.Range("E19").Value = '  +0.02%  '

$ws.Range("B20").Value = 'Litecoin'
$ws.Range("C20").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '81.32'
$c.Style = "Normal"
$ws.Range("E20").Value = '  +4.29%  '

$ws.Range("B21").Value = 'Avalanche'
$ws.Range("C21").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '16.46'
$c.Style = "Normal"
$ws.Range("E21").Value = '  +2.20%  '

$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '6.112'
$c.Style = "Normal"
$ws.Range("E22").Value = '  +1.51%  '

$ws.Range("B23").Value = 'Cosmos'
$ws.Range("C23").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '11.99'
$c.Style = "Normal"
$ws.Range("E23").Value = '  +1.50%  '

$ws.Range("B24").Value = 'WrappedBTC'
$ws.Range("C24").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D24").Value = '24.394.75'
$ws.Range("E24").Value = '  +8.66%  '

$ws.Range("B25").Value = 'Toncoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '2.445'
$c.Style = "Normal"
$ws.Range("E25").Value = '  +1.83%  '

$ws.Range("B26").Value = 'LEO'
$ws.Range("C26").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '3.361'
$c.Style = "Normal"
$ws.Range("E26").Value = '  -11.72%  '

$ws.Range("B27").Value = 'LidoDAOToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '2.675'
$c.Style = "Normal"
$ws.Range("E27").Value = '  +6.57%  '

$ws.Range("B28").Value = 'Monero'
$ws.Range("C28").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '152.56'
$c.Style = "Normal"
$ws.Range("E28").Value = '  +1.29%  '

$ws.Range("B29").Value = 'EthereumClassic'
$ws.Range("C29").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '19.59'
$c.Style = "Normal"
$ws.Range("E29").Value = '  +0.36%  '

$ws.Range("B30").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C30").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D30").Value = '1.859.15'
$ws.Range("E30").Value = '  +4.32%  '

$ws.Range("B31").Value = 'BitcoinCash'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '126.93'
$c.Style = "Normal"
$ws.Range("E31").Value = '  +5.18%  '

$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '6.359'
$c.Style = "Normal"
$ws.Range("E32").Value = '  +5.38%  '

$ws.Range("B33").Value = 'HuobiToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '4.038'
$c.Style = "Normal"
$ws.Range("E33").Value = '  -3.81%  '

$ws.Range("B34").Value = 'ImmutableX'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '0.9736'
$c.Style = "Normal"
$ws.Range("E34").Value = '  +1.78%  '

$ws.Range("B35").Value = 'Stellar'
$ws.Range("C35").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '0.08480'
$c.Style = "Normal"
$ws.Range("E35").Value = '  +2.45%  '

$ws.Range("B36").Value = 'WEMIXTOKEN'
$ws.Range("C36").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '1.673'
$c.Style = "Normal"
$ws.Range("E36").Value = '  +2.20%  '

$ws.Range("B37").Value = 'Aptos'
$ws.Range("C37").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '12.47'
$c.Style = "Normal"
$ws.Range("E37").Value = '  +5.21%  '

$ws.Range("B38").Value = 'Hedera'
$ws.Range("C38").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.06498'
$c.Style = "Normal"
$ws.Range("E38").Value = '  +6.36%  '

$ws.Range("B39").Value = 'FraxShare'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '8.958'
$c.Style = "Normal"
$ws.Range("E39").Value = '  +4.23%  '

$ws.Range("B40").Value = 'InternetComputer(DFINITY)'
$ws.Range("C40").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '5.359'
$c.Style = "Normal"
$ws.Range("E40").Value = '  +1.31%  '

$ws.Range("B41").Value = 'VeChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '0.02339'
$c.Style = "Normal"
$ws.Range("E41").Value = '  +5.40%  '

$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '1.272'
$c.Style = "Normal"
$ws.Range("E42").Value = '  -0.10%  '

$ws.Range("B43").Value = 'Algorand'
$ws.Range("C43").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '0.2116'
$c.Style = "Normal"
$ws.Range("E43").Value = '  +4.32%  '

$ws.Range("B44").Value = 'TheSandbox'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '0.6189'
$c.Style = "Normal"
$ws.Range("E44").Value = '  +4.68%  '

$ws.Range("B45").Value = 'Frax'
$ws.Range("C45").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '0.9982'
$c.Style = "Normal"
$ws.Range("E45").Value = '  +0.68%  '

$ws.Range("B46").Value = 'PancakeSwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '3.771'
$c.Style = "Normal"
$ws.Range("E46").Value = '  -2.27%  '

$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '13.09'
$c.Style = "Normal"
$ws.Range("E47").Value = '  -1.35%  '

$ws.Range("B48").Value = 'Decentraland'
$ws.Range("C48").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '0.5956'
$c.Style = "Normal"
$ws.Range("E48").Value = '  +4.63%  '

$ws.Range("B49").Value = 'Quant'
$ws.Range("C49").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '127.11'
$c.Style = "Normal"
$ws.Range("E49").Value = '  +0.06%  '

$ws.Range("B50").Value = 'NEARProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '2.034'
$c.Style = "Normal"
$ws.Range("E50").Value = '  +3.32%  '

$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '0.07219'
$c.Style = "Normal"
$ws.Range("E51").Value = '  +5.79%  '
